$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update environment-specific values from test7 -> test14
$ws.Range("A2").Value = "https://test14.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test14.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test14.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest14"
$ws.Range("G2").Value = "test14"
$ws.Range("K2").Value = "virtual_cabitest14"

# Update the active cell selection on the sheet view
$ws.Range("C11").Select()
